$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-observation row is inserted above row 543, shifting the
# existing rows 543-593 down to 544-594 (row 594 ends up holding the data
# that used to be in row 593).
$ws.Rows.Item(543).Insert()

# Populate the newly inserted row 543 with the new observation. All columns
# besides the price/date fields (A, B, C, E, F, G, H, I, N, O, Q, R) follow
# the same constant pattern as every other row in this block.
$ws.Range("A543").Value = 8
$ws.Range("B543").Value = "Terminal La Palmera de La Serena"
$ws.Range("C543").Value = "Coquimbo"
$ws.Range("D543").Value = 45166
$ws.Range("E543").Value = 4
$ws.Range("F543").Value = 100114013
$ws.Range("G543").Value = "Zanahoria"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 440
$ws.Range("K543").Value = 5500
$ws.Range("L543").Value = 6000
$ws.Range("M543").Value = 5750
$ws.Range("N543").Value = "$/saco 20 kilos"
$ws.Range("O543").Value = "Provincia del Elquí"
$ws.Range("P543").Value = 288
$ws.Range("Q543").Value = 20
$ws.Range("R543").Value = "Hortaliza"
